# Update the crypto price/volume table with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $savedStyle = $rng.Style
    $rng.Value = "'" + $value
    $rng.Style = $savedStyle
}

# Row 2
Set-TextValue "D2" "65.855.69"
Set-TextValue "E2" "  +0.34%  "

# Row 3
Set-TextValue "D3" "2.663.74"
Set-TextValue "E3" "  -0.43%  "

# Row 4
Set-TextValue "E4" "  -0.01%  "

# Row 5
Set-TextValue "D5" "598.39"
Set-TextValue "E5" "  -0.36%  "

# Row 6
Set-TextValue "D6" "157.82"
Set-TextValue "E6" "  +0.74%  "

# Row 7
Set-TextValue "E7" "  +4.61%  "

# Row 8
Set-TextValue "E8" "  -0.02%  "

# Row 9
Set-TextValue "E9" "  -1.90%  "

# Row 10
Set-TextValue "D10" "0.403"
Set-TextValue "E10" "  +0.48%  "

# Row 11
Set-TextValue "E11" "  +0.03%  "

# Row 12
Set-TextValue "E12" "  +1.64%  "

# Row 13
Set-TextValue "D13" "29.05"
Set-TextValue "E13" "  -1.05%  "

# Row 14
Set-TextValue "E14" "  -1.77%  "

# Row 15
Set-TextValue "D15" "3.143.64"
Set-TextValue "E15" "  -0.36%  "

# Row 16
Set-TextValue "D16" "65.712.12"
Set-TextValue "E16" "  +0.34%  "

# Row 17
Set-TextValue "D17" "2.658.99"
Set-TextValue "E17" "  -0.40%  "

# Row 18
Set-TextValue "D18" "12.64"
Set-TextValue "E18" "  -1.37%  "

# Rows 20/21 - BitcoinCash and Uniswap swap positions with new values
Set-TextValue "B20" "Uniswap"
Set-TextValue "C20" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D20" "7.50"
Set-TextValue "E20" "  -0.93%  "

Set-TextValue "B21" "BitcoinCash"
Set-TextValue "C21" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D21" "351.61"
Set-TextValue "E21" "  -0.06%  "

# Row 23
Set-TextValue "D23" "69.83"

# Row 24
Set-TextValue "D24" "1.84"
Set-TextValue "E24" "  +12.19%  "

# Row 25
Set-TextValue "E25" "  +0.90%  "

# Row 26
Set-TextValue "D26" "9.74"
Set-TextValue "E26" "  +0.90%  "

# Row 27
Set-TextValue "D27" "1.63"
Set-TextValue "E27" "  +1.50%  "

# Row 28
Set-TextValue "D28" "574.36"
Set-TextValue "E28" "  +8.12%  "

# Row 29
Set-TextValue "D29" "8.19"
Set-TextValue "E29" "  +1.12%  "

# Row 30
Set-TextValue "D30" "0.163"
Set-TextValue "E30" "  -2.36%  "

# Row 31
Set-TextValue "D31" "0.999"
Set-TextValue "E31" "  -0.18%  "

# Row 32
Set-TextValue "E32" "  -0.26%  "

# Row 33
Set-TextValue "D33" "1.83"
Set-TextValue "E33" "  +4.26%  "

# Row 34
Set-TextValue "D34" "6.67"
Set-TextValue "E34" "  +3.12%  "

# Row 35
Set-TextValue "D35" "5.61"
Set-TextValue "E35" "  +2.08%  "

# Row 36
Set-TextValue "E36" "  -0.28%  "

# Row 37
Set-TextValue "D37" "20.56"
Set-TextValue "E37" "  +0.12%  "

# Row 38
Set-TextValue "E38" "  -0.01%  "

# Row 39
Set-TextValue "E39" "  +0.79%  "

# Row 40
Set-TextValue "D40" "154.50"
Set-TextValue "E40" "  -2.14%  "

# Row 41
Set-TextValue "D41" "161.16"
Set-TextValue "E41" "  -1.93%  "

# Row 42
Set-TextValue "E42" "  -1.68%  "

# Row 43
Set-TextValue "D43" "0.0617"
Set-TextValue "E43" "  +1.30%  "

# Row 44
Set-TextValue "E44" "  +0.23%  "

# Row 45
Set-TextValue "D45" "23.03"
Set-TextValue "E45" "  +0.78%  "

# Row 46
Set-TextValue "E46" "  +0.23%  "

# Row 47
Set-TextValue "D47" "0.0258"
Set-TextValue "E47" "  -0.76%  "

# Row 48
Set-TextValue "E48" "  +0.90%  "

# Row 49
Set-TextValue "D49" "19.82"
Set-TextValue "E49" "  -1.57%  "

# Row 50
Set-TextValue "E50" "  -5.91%  "

# Row 51
Set-TextValue "D51" "0.820"
Set-TextValue "E51" "  +0.40%  "
